$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 288.77777
$ws.Range("I6").Value = 157
$ws.Range("K6").Value = 471
$ws.Range("M6").Value = -359
$ws.Range("H92").Value = 236.89473
$ws.Range("I92").Value = 212.11765
$ws.Range("J92").Value = 447.5
$ws.Range("K92").Value = 212.11765
$ws.Range("L92").Value = 447.5
$ws.Range("M92").Value = 1035.88235
$ws.Range("N92").Value = -2943.5
$ws.Range("H133").Value = 50059.375
$ws.Range("J133").Value = 50059.375
$ws.Range("L133").Value = 50059.375
$ws.Range("N133").Value = -60179.375
$ws.Range("H138").Value = 6226.1113
$ws.Range("I138").Value = 4209.6
$ws.Range("J138").Value = 6452.6855
$ws.Range("K138").Value = 12628.8
$ws.Range("L138").Value = 19358.0565
$ws.Range("M138").Value = -7488.800000000001
$ws.Range("N138").Value = -29638.0565
$ws.Range("H139").Value = 39756.855
$ws.Range("J139").Value = 39756.855
$ws.Range("L139").Value = 39756.855
$ws.Range("N139").Value = -50036.855

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1634.5172
$ws.Range("I132").Value = 1466.7407
$ws.Range("J132").Value = 3899.5
$ws.Range("K132").Value = 4400.2221
$ws.Range("L132").Value = 11698.5
$ws.Range("M132").Value = -1870.2221
$ws.Range("N132").Value = -16758.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 42500
$ws.Range("J122").Value = 42500
$ws.Range("L122").Value = 42500
$ws.Range("N122").Value = -52300
$ws.Range("H126").Value = 27955
$ws.Range("J126").Value = 27955
$ws.Range("L126").Value = 27955
$ws.Range("N126").Value = -37835

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 68093
$ws.Range("I31").Value = 113957.11
$ws.Range("J31").Value = 9124.857
$ws.Range("K31").Value = 113957.11
$ws.Range("L31").Value = 9124.857
$ws.Range("M31").Value = -113662.11
$ws.Range("N31").Value = -9714.857
$ws.Range("H34").Value = 68093
$ws.Range("I34").Value = 113957.11
$ws.Range("J34").Value = 9124.857
$ws.Range("K34").Value = 113957.11
$ws.Range("L34").Value = 9124.857
$ws.Range("M34").Value = -113755.11
$ws.Range("N34").Value = -9528.857
$ws.Range("H50").Value = 7851.143
$ws.Range("J50").Value = 9145.833
$ws.Range("L50").Value = 9145.833
$ws.Range("N50").Value = -10395.833
$ws.Range("H68").Value = 15026.75
$ws.Range("J68").Value = 17491
$ws.Range("L68").Value = 17491
$ws.Range("N68").Value = -18989
$ws.Range("H71").Value = 15026.75
$ws.Range("J71").Value = 17491
$ws.Range("L71").Value = 52473
$ws.Range("N71").Value = -59961
$ws.Range("H130").Value = 56086.668
$ws.Range("J130").Value = 56086.668
$ws.Range("L130").Value = 56086.668
$ws.Range("N130").Value = -66126.668

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 277.14285
$ws.Range("I7").Value = 163.33333
$ws.Range("J7").Value = 362.5
$ws.Range("K7").Value = 489.99999
$ws.Range("L7").Value = 1087.5
$ws.Range("M7").Value = -377.99999
$ws.Range("N7").Value = -1311.5
$ws.Range("H17").Value = 11400.111
$ws.Range("J17").Value = 15669.385
$ws.Range("L17").Value = 47008.155
$ws.Range("N17").Value = -47346.155
$ws.Range("H34").Value = 883.6087
$ws.Range("I34").Value = 316.625
$ws.Range("J34").Value = 1186
$ws.Range("K34").Value = 949.875
$ws.Range("L34").Value = 3558
$ws.Range("M34").Value = -865.875
$ws.Range("N34").Value = -3726
$ws.Range("H39").Value = 2242.2666
$ws.Range("J39").Value = 2527.2307
$ws.Range("L39").Value = 7581.6921
$ws.Range("N39").Value = -8169.6921
$ws.Range("H55").Value = 93009.55
$ws.Range("J55").Value = 102300.5
$ws.Range("L55").Value = 306901.5
$ws.Range("N55").Value = -307255.5
$ws.Range("H75").Value = 2755.88
$ws.Range("J75").Value = 2812.375
$ws.Range("L75").Value = 8437.125
$ws.Range("N75").Value = -10433.125
$ws.Range("H78").Value = 2755.88
$ws.Range("J78").Value = 2812.375
$ws.Range("L78").Value = 25311.375
$ws.Range("N78").Value = -35295.375
$ws.Range("H80").Value = 8191.8335
$ws.Range("I80").Value = 5667.3335
$ws.Range("J80").Value = 9033.333
$ws.Range("K80").Value = 17002.0005
$ws.Range("L80").Value = 27099.999
$ws.Range("M80").Value = -16066.0005
$ws.Range("N80").Value = -28971.999
$ws.Range("H83").Value = 8191.8335
$ws.Range("I83").Value = 5667.3335
$ws.Range("J83").Value = 9033.333
$ws.Range("K83").Value = 51006.0015
$ws.Range("L83").Value = 81299.997
$ws.Range("M83").Value = -46326.0015
$ws.Range("N83").Value = -90659.997
$ws.Range("H92").Value = 3333.3333
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 3333.3333
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 9999.999899999999
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -12495.9999
$ws.Range("H122").Value = 1520.7778
$ws.Range("I122").Value = 760.4
$ws.Range("J122").Value = 1813.2307
$ws.Range("K122").Value = 6843.599999999999
$ws.Range("L122").Value = 16319.0763
$ws.Range("M122").Value = -4393.599999999999
$ws.Range("N122").Value = -21219.0763
$ws.Range("H131").Value = 847.98
$ws.Range("I131").Value = 540.9167
$ws.Range("J131").Value = 889.8523
$ws.Range("K131").Value = 1622.7501
$ws.Range("L131").Value = 2669.5569
$ws.Range("M131").Value = 3417.2499
$ws.Range("N131").Value = -12749.5569

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 21173.666
$ws.Range("J57").Value = 20060.5
$ws.Range("L57").Value = 20060.5
$ws.Range("N57").Value = -21700.5
$ws.Range("H132").Value = 2775.6128
$ws.Range("I132").Value = 2722.1667
$ws.Range("J132").Value = 2958.8572
$ws.Range("K132").Value = 8166.500100000001
$ws.Range("L132").Value = 8876.5716
$ws.Range("M132").Value = -5636.500100000001
$ws.Range("N132").Value = -13936.5716

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3191.32
$ws.Range("I132").Value = 2938.9333
$ws.Range("J132").Value = 3569.9
$ws.Range("K132").Value = 8816.7999
$ws.Range("L132").Value = 10709.7
$ws.Range("M132").Value = -6286.7999
$ws.Range("N132").Value = -15769.7

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 22377
$ws.Range("J109").Value = 22377
$ws.Range("L109").Value = 22377
$ws.Range("N109").Value = -25151
$ws.Range("H123").Value = 50714.5
$ws.Range("J123").Value = 50714.5
$ws.Range("L123").Value = 50714.5
$ws.Range("N123").Value = -60514.5
$ws.Range("H132").Value = 1882.3405
$ws.Range("I132").Value = 1802.6666
$ws.Range("J132").Value = 2022.9412
$ws.Range("K132").Value = 5407.9998
$ws.Range("L132").Value = 6068.8236
$ws.Range("M132").Value = -2877.9998
$ws.Range("N132").Value = -11128.8236
